$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 11).Value = 13

# Row 8
$ws.Cells.Item(8, 10).Value = 1.04
$ws.Cells.Item(8, 11).Value = 13
$ws.Cells.Item(8, 14).Value = 1.75
$ws.Cells.Item(8, 15).Value = 2.05

# Row 9
$ws.Cells.Item(9, 7).Value = 2.35
$ws.Cells.Item(9, 8).Value = 2.6
$ws.Cells.Item(9, 9).Value = 3.75
$ws.Cells.Item(9, 11).Value = 5.2
$ws.Cells.Item(9, 12).Value = 1.47
$ws.Cells.Item(9, 13).Value = 2.5
$ws.Cells.Item(9, 16).Value = 1.53
$ws.Cells.Item(9, 17).Value = 2.35
$ws.Cells.Item(9, 21).Value = 10.75
$ws.Cells.Item(9, 22).Value = 9
$ws.Cells.Item(9, 23).Value = 26
$ws.Cells.Item(9, 24).Value = 22
$ws.Cells.Item(9, 25).Value = 35
$ws.Cells.Item(9, 26).Value = 5.2
$ws.Cells.Item(9, 31).Value = 8.5
$ws.Cells.Item(9, 32).Value = 20
$ws.Cells.Item(9, 33).Value = 12.5
$ws.Cells.Item(9, 34).Value = 65

# Row 10
$ws.Cells.Item(10, 7).Value = 1.8
$ws.Cells.Item(10, 8).Value = 3.3
$ws.Cells.Item(10, 9).Value = 4.55
$ws.Cells.Item(10, 10).Value = 1.1
$ws.Cells.Item(10, 11).Value = 6
$ws.Cells.Item(10, 12).Value = 1.45
$ws.Cells.Item(10, 13).Value = 2.57
$ws.Cells.Item(10, 14).Value = 2.3
$ws.Cells.Item(10, 15).Value = 1.55
$ws.Cells.Item(10, 16).Value = 1.5
$ws.Cells.Item(10, 17).Value = 2.42
$ws.Cells.Item(10, 18).Value = 2.15
$ws.Cells.Item(10, 19).Value = 1.62
$ws.Cells.Item(10, 20).Value = 5.3
$ws.Cells.Item(10, 21).Value = 7.2
$ws.Cells.Item(10, 22).Value = 8.75
$ws.Cells.Item(10, 23).Value = 14
$ws.Cells.Item(10, 24).Value = 17
$ws.Cells.Item(10, 25).Value = 40
$ws.Cells.Item(10, 26).Value = 6
$ws.Cells.Item(10, 27).Value = 6.5
$ws.Cells.Item(10, 28).Value = 19.5
$ws.Cells.Item(10, 29).Value = 120
$ws.Cells.Item(10, 31).Value = 10
$ws.Cells.Item(10, 32).Value = 23
$ws.Cells.Item(10, 33).Value = 15.5
$ws.Cells.Item(10, 34).Value = 80
$ws.Cells.Item(10, 35).Value = 55
$ws.Cells.Item(10, 36).Value = 65

# Row 12
$ws.Cells.Item(12, 7).Value = 3.25
$ws.Cells.Item(12, 9).Value = 2.05
$ws.Cells.Item(12, 10).Value = 1.06
$ws.Cells.Item(12, 11).Value = 10
$ws.Cells.Item(12, 23).Value = 41
$ws.Cells.Item(12, 31).Value = 6.5
$ws.Cells.Item(12, 33).Value = 9

# Row 15
$ws.Cells.Item(15, 11).Value = 10
$ws.Cells.Item(15, 12).Value = 1.3
$ws.Cells.Item(15, 13).Value = 3.4
$ws.Cells.Item(15, 14).Value = 2
$ws.Cells.Item(15, 15).Value = 1.8

# Row 17
$ws.Cells.Item(17, 7).Value = 3.8
$ws.Cells.Item(17, 8).Value = 3.25
$ws.Cells.Item(17, 9).Value = 1.93
$ws.Cells.Item(17, 12).Value = 1.36
$ws.Cells.Item(17, 13).Value = 2.67
$ws.Cells.Item(17, 14).Value = 2.05
$ws.Cells.Item(17, 15).Value = 1.6
$ws.Cells.Item(17, 16).Value = 1.42
$ws.Cells.Item(17, 17).Value = 2.47
$ws.Cells.Item(17, 18).Value = 1.87
$ws.Cells.Item(17, 19).Value = 1.75
$ws.Cells.Item(17, 20).Value = 9.25
$ws.Cells.Item(17, 21).Value = 19.5
$ws.Cells.Item(17, 22).Value = 13.5
$ws.Cells.Item(17, 23).Value = 60
$ws.Cells.Item(17, 24).Value = 40
$ws.Cells.Item(17, 25).Value = 50
$ws.Cells.Item(17, 26).Value = 8.25
$ws.Cells.Item(17, 27).Value = 6.3
$ws.Cells.Item(17, 28).Value = 16.5
$ws.Cells.Item(17, 29).Value = 90
$ws.Cells.Item(17, 30).Value = 800
$ws.Cells.Item(17, 31).Value = 6.3
$ws.Cells.Item(17, 32).Value = 8.5
$ws.Cells.Item(17, 33).Value = 8.5
$ws.Cells.Item(17, 34).Value = 16.5
$ws.Cells.Item(17, 35).Value = 16.5
$ws.Cells.Item(17, 36).Value = 32

# Row 18
$ws.Cells.Item(18, 7).Value = 5.4
$ws.Cells.Item(18, 8).Value = 3.7
$ws.Cells.Item(18, 9).Value = 1.57
$ws.Cells.Item(18, 12).Value = 1.25
$ws.Cells.Item(18, 13).Value = 3.2
$ws.Cells.Item(18, 14).Value = 1.75
$ws.Cells.Item(18, 15).Value = 1.85
$ws.Cells.Item(18, 18).Value = 1.78
$ws.Cells.Item(18, 19).Value = 1.83
$ws.Cells.Item(18, 20).Value = 14.5
$ws.Cells.Item(18, 21).Value = 32
$ws.Cells.Item(18, 22).Value = 17
$ws.Cells.Item(18, 23).Value = 100
$ws.Cells.Item(18, 24).Value = 55
$ws.Cells.Item(18, 25).Value = 55
$ws.Cells.Item(18, 26).Value = 10.5
$ws.Cells.Item(18, 27).Value = 7.3
$ws.Cells.Item(18, 28).Value = 16
$ws.Cells.Item(18, 29).Value = 75
$ws.Cells.Item(18, 30).Value = 600
$ws.Cells.Item(18, 31).Value = 6.9
$ws.Cells.Item(18, 32).Value = 7.5
$ws.Cells.Item(18, 33).Value = 8
$ws.Cells.Item(18, 34).Value = 11.75
$ws.Cells.Item(18, 35).Value = 12.5
$ws.Cells.Item(18, 36).Value = 25

# Row 19
$ws.Cells.Item(19, 7).Value = 3.5
$ws.Cells.Item(19, 8).Value = 2.77
$ws.Cells.Item(19, 9).Value = 2.25
$ws.Cells.Item(19, 12).Value = 1.44
$ws.Cells.Item(19, 13).Value = 2.42
$ws.Cells.Item(19, 14).Value = 2.25
$ws.Cells.Item(19, 15).Value = 1.5
$ws.Cells.Item(19, 16).Value = 1.55
$ws.Cells.Item(19, 17).Value = 2.15
$ws.Cells.Item(19, 18).Value = 1.88
$ws.Cells.Item(19, 19).Value = 1.72
$ws.Cells.Item(19, 20).Value = 8.75
$ws.Cells.Item(19, 21).Value = 18.5
$ws.Cells.Item(19, 22).Value = 12
$ws.Cells.Item(19, 23).Value = 55
$ws.Cells.Item(19, 24).Value = 35
$ws.Cells.Item(19, 25).Value = 45
$ws.Cells.Item(19, 26).Value = 6.6
$ws.Cells.Item(19, 27).Value = 5.5
$ws.Cells.Item(19, 28).Value = 15
$ws.Cells.Item(19, 29).Value = 90
$ws.Cells.Item(19, 30).Value = 800
$ws.Cells.Item(19, 31).Value = 6
$ws.Cells.Item(19, 32).Value = 9.75
$ws.Cells.Item(19, 33).Value = 9.25
$ws.Cells.Item(19, 34).Value = 23
$ws.Cells.Item(19, 35).Value = 22
$ws.Cells.Item(19, 36).Value = 37

# Row 24
$ws.Cells.Item(24, 14).Value = 2.1
$ws.Cells.Item(24, 15).Value = 1.7

# Row 26
$ws.Cells.Item(26, 11).Value = 9

# Row 27
$ws.Cells.Item(27, 7).Value = 1.88
$ws.Cells.Item(27, 8).Value = 3.3
$ws.Cells.Item(27, 9).Value = 3.95
$ws.Cells.Item(27, 12).Value = 1.44
$ws.Cells.Item(27, 13).Value = 2.42
$ws.Cells.Item(27, 14).Value = 2.25
$ws.Cells.Item(27, 15).Value = 1.5
$ws.Cells.Item(27, 16).Value = 1.5
$ws.Cells.Item(27, 17).Value = 2.25
$ws.Cells.Item(27, 18).Value = 2.07
$ws.Cells.Item(27, 19).Value = 1.6
$ws.Cells.Item(27, 20).Value = 5.5
$ws.Cells.Item(27, 21).Value = 7.7
$ws.Cells.Item(27, 22).Value = 9
$ws.Cells.Item(27, 23).Value = 15.5
$ws.Cells.Item(27, 24).Value = 18
$ws.Cells.Item(27, 25).Value = 40
$ws.Cells.Item(27, 26).Value = 7.3
$ws.Cells.Item(27, 27).Value = 6.5
$ws.Cells.Item(27, 28).Value = 20
$ws.Cells.Item(27, 29).Value = 120
$ws.Cells.Item(27, 31).Value = 8.75
$ws.Cells.Item(27, 32).Value = 19.5
$ws.Cells.Item(27, 33).Value = 14.5
$ws.Cells.Item(27, 34).Value = 60
$ws.Cells.Item(27, 35).Value = 45
$ws.Cells.Item(27, 36).Value = 65

# Row 28
$ws.Cells.Item(28, 7).Value = 2.95
$ws.Cells.Item(28, 8).Value = 2.95
$ws.Cells.Item(28, 9).Value = 2.45
$ws.Cells.Item(28, 12).Value = 1.52
$ws.Cells.Item(28, 13).Value = 2.22
$ws.Cells.Item(28, 14).Value = 2.45
$ws.Cells.Item(28, 15).Value = 1.42
$ws.Cells.Item(28, 16).Value = 1.57
$ws.Cells.Item(28, 17).Value = 2.1
$ws.Cells.Item(28, 18).Value = 2.07
$ws.Cells.Item(28, 19).Value = 1.6
$ws.Cells.Item(28, 20).Value = 6.8
$ws.Cells.Item(28, 21).Value = 13
$ws.Cells.Item(28, 22).Value = 11.5
$ws.Cells.Item(28, 23).Value = 37
$ws.Cells.Item(28, 24).Value = 32
$ws.Cells.Item(28, 25).Value = 55
$ws.Cells.Item(28, 26).Value = 6.3
$ws.Cells.Item(28, 27).Value = 5.9
$ws.Cells.Item(28, 28).Value = 19
$ws.Cells.Item(28, 29).Value = 120
$ws.Cells.Item(28, 31).Value = 5.9
$ws.Cells.Item(28, 32).Value = 10.25
$ws.Cells.Item(28, 33).Value = 10.25
$ws.Cells.Item(28, 34).Value = 26
$ws.Cells.Item(28, 35).Value = 26
$ws.Cells.Item(28, 36).Value = 50

# Row 30
$ws.Cells.Item(30, 10).Value = 1.01
$ws.Cells.Item(30, 11).Value = 13
$ws.Cells.Item(30, 18).Value = 1.57
$ws.Cells.Item(30, 19).Value = 2.25
$ws.Cells.Item(30, 20).Value = 13
$ws.Cells.Item(30, 21).Value = 19
$ws.Cells.Item(30, 28).Value = 12
$ws.Cells.Item(30, 30).Value = 126
$ws.Cells.Item(30, 31).Value = 10
$ws.Cells.Item(30, 32).Value = 12

# Row 33
$ws.Cells.Item(33, 7).Value = 1.62
$ws.Cells.Item(33, 8).Value = 4.33
$ws.Cells.Item(33, 10).Value = 21
$ws.Cells.Item(33, 11).Value = 1.03
$ws.Cells.Item(33, 12).Value = 1.13
$ws.Cells.Item(33, 13).Value = 5.5
$ws.Cells.Item(33, 20).Value = 12
$ws.Cells.Item(33, 27).Value = 9.5
$ws.Cells.Item(33, 31).Value = 21

# Row 36
$ws.Cells.Item(36, 7).Value = 2.63
$ws.Cells.Item(36, 9).Value = 2.3
$ws.Cells.Item(36, 16).Value = 1.3
$ws.Cells.Item(36, 17).Value = 3.4
$ws.Cells.Item(36, 24).Value = 21
$ws.Cells.Item(36, 33).Value = 9.5
$ws.Cells.Item(36, 35).Value = 17

# Row 38
$ws.Cells.Item(38, 10).Value = 1.02
$ws.Cells.Item(38, 11).Value = 12

# Row 39
$ws.Cells.Item(39, 14).Value = 2.05
$ws.Cells.Item(39, 15).Value = 1.75

# Row 40
$ws.Cells.Item(40, 7).Value = 3.3
$ws.Cells.Item(40, 8).Value = 3.75
$ws.Cells.Item(40, 9).Value = 2.05
$ws.Cells.Item(40, 18).Value = 1.57
$ws.Cells.Item(40, 19).Value = 2.25
$ws.Cells.Item(40, 20).Value = 13
$ws.Cells.Item(40, 22).Value = 12
$ws.Cells.Item(40, 25).Value = 29
$ws.Cells.Item(40, 32).Value = 11

# Row 41
$ws.Cells.Item(41, 10).Value = 1.03
$ws.Cells.Item(41, 11).Value = 15

# Row 46
$ws.Cells.Item(46, 12).Value = 1.06
$ws.Cells.Item(46, 13).Value = 10
